$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (employee/training id 6) had a placeholder text value ("Completo")
# in the "validade" (D) column. Replace it with an actual expiration date,
# matching the date format already used by the other rows/columns
# (numFmt "yyyy-mm-dd", same as column C).
$ws.Range("D7").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("D7").Value = 45087

# Reflect where the user's cursor ended up after making the edit.
$ws.Range("B11").Select()
